$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0007018076591000852
$ws.Range("E2").Value = 0.0007018076591000852
$ws.Range("D3").Value = 0.9581194954748646
$ws.Range("E3").Value = 0.9581194954748646
$ws.Range("D4").Value = 0.00008540584901100519
$ws.Range("E4").Value = 0.00008540584901100519
$ws.Range("D5").Value = 0.0000001367006494339108
$ws.Range("E5").Value = 0.0000001367006494339108
$ws.Range("D6").Value = 0.1784996208733365
$ws.Range("E6").Value = 0.1784996208733365
$ws.Range("D7").Value = 0.9728231196496859
$ws.Range("E7").Value = 0.02717688035031407
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.0001245355230851506
$ws.Range("E8").Value = 0.9998754644769149
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.00006688183016069859
$ws.Range("E9").Value = 0.9999331181698393
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.103599041923544
$ws.Range("E10").Value = 0.896400958076456
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.00003675190955152106
$ws.Range("E11").Value = 0.9999632480904485
$ws.Range("F11").Value = 3.447994947433472
$ws.Range("G11").Value = 0.5
$ws.Range("D12").Value = 0.000000304118327605434
$ws.Range("E12").Value = 0.000000304118327605434
$ws.Range("D13").Value = 0.9953023132390378
$ws.Range("E13").Value = 0.9953023132390378
$ws.Range("D14").Value = 0.00003708685082327059
$ws.Range("E14").Value = 0.00003708685082327059
$ws.Range("D15").Value = 0.0000000001076239516145161
$ws.Range("E15").Value = 0.0000000001076239516145161
$ws.Range("D16").Value = 0.09142299564742942
$ws.Range("E16").Value = 0.09142299564742942
$ws.Range("D17").Value = 0.9712806692736689
$ws.Range("E17").Value = 0.02871933072633115
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.0000001694043086397785
$ws.Range("E18").Value = 0.9999998305956913
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.000006290077407956476
$ws.Range("E19").Value = 0.9999937099225921
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.3351822036917832
$ws.Range("E20").Value = 0.6648177963082168
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.00000004051695403488507
$ws.Range("E21").Value = 0.999999959483046
$ws.Range("F21").Value = 5.116787910461426
$ws.Range("G21").Value = 0.5
